# Articulation Test Refactoring
# - Position column (C) now holds mixed-case values ("Initial", "Medial",
#   "Final", "Blended") instead of the old all-caps strings.
# - Image column (D) is reshuffled between "lemon.jpg" and "book.jpg"
#   (replacing the single "lemon.png" placeholder).
# - The generated SQL (column E / the CONCATENATE formula) now targets the
#   renamed [SoundPosition] column instead of [Position], and the sample
#   text reflects the new Position/Image values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- Column C: Position values, mixed case, per block of rows ---
$ws.Range("C2:C21").Value  = "Initial"
$ws.Range("C22:C43").Value = "Medial"
$ws.Range("C44:C63").Value = "Final"
$ws.Range("C64:C81").Value = "Blended"

# --- Column D: Image filename, set individually per row ---
$images = @(
    "book.jpg","lemon.jpg","lemon.jpg","book.jpg","lemon.jpg","book.jpg","lemon.jpg","lemon.jpg","lemon.jpg","lemon.jpg",
    "book.jpg","lemon.jpg","book.jpg","lemon.jpg","lemon.jpg","lemon.jpg","lemon.jpg","lemon.jpg","book.jpg","book.jpg",
    "lemon.jpg","book.jpg","lemon.jpg","lemon.jpg","lemon.jpg","lemon.jpg","lemon.jpg","lemon.jpg","lemon.jpg","lemon.jpg",
    "book.jpg","lemon.jpg","book.jpg","lemon.jpg","lemon.jpg","lemon.jpg","lemon.jpg","book.jpg","lemon.jpg","lemon.jpg",
    "book.jpg","lemon.jpg","lemon.jpg","lemon.jpg","lemon.jpg","lemon.jpg","lemon.jpg","book.jpg","lemon.jpg","lemon.jpg",
    "lemon.jpg","lemon.jpg","lemon.jpg","book.jpg","lemon.jpg","lemon.jpg","book.jpg","lemon.jpg","lemon.jpg","lemon.jpg",
    "lemon.jpg","lemon.jpg","lemon.jpg","lemon.jpg","lemon.jpg","book.jpg","lemon.jpg","lemon.jpg","book.jpg","lemon.jpg",
    "lemon.jpg","lemon.jpg","lemon.jpg","book.jpg","lemon.jpg","lemon.jpg","book.jpg","lemon.jpg","lemon.jpg","book.jpg"
)

for ($i = 0; $i -lt $images.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $images[$i]
}

# --- Column E: formula text, [Position] -> [SoundPosition] ---
# Row 2 has its own (non-shared) formula.
$ws.Range("E2").Formula = "=CONCATENATE(""INSERT INTO [ArticulationTests] ([Sound],[Text],[SoundPosition],[Image])  VALUES ('"",A2,""','"",B2,""','"",C2,""','"",D2,""');"")"

# Rows 3:66 share one formula (master at E3) - set the whole block together
# so it stays a single shared formula, matching the original layout.
$ws.Range("E3:E66").Formula = "=CONCATENATE(""INSERT INTO [ArticulationTests] ([Sound],[Text],[SoundPosition],[Image])  VALUES ('"",A3,""','"",B3,""','"",C3,""','"",D3,""');"")"

# Rows 67:81 share a second formula (master at E67).
$ws.Range("E67:E81").Formula = "=CONCATENATE(""INSERT INTO [ArticulationTests] ([Sound],[Text],[SoundPosition],[Image])  VALUES ('"",A67,""','"",B67,""','"",C67,""','"",D67,""');"")"
